$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 82925000
$ws.Range("D3").Value = 47800000
$ws.Range("D4").Value = 40300000
$ws.Range("D5").Value = 32800000

$ws.Range("B6").Value = 0.5000000000039895
$ws.Range("C6").Value = 0.5000000000039895
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
